$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestSteps")

# Row 7 (A:C) previously had no data/style in those columns; give it the
# same border/fill formatting already used by the other keyword rows
# before we populate it below.
$ws.Range("A6:C6").Copy() | Out-Null
$ws.Range("A7:C7").PasteSpecial(-4122) | Out-Null

# Shift the existing keyword rows down by one (rows 3-6 -> rows 4-7),
# working from the bottom up so we don't clobber source data.
$ws.Range("A7").Value = $ws.Range("A6").Value()
$ws.Range("B7").Value = $ws.Range("B6").Value()
$ws.Range("C7").Value = $ws.Range("C6").Value()

$ws.Range("A6").Value = $ws.Range("A5").Value()
$ws.Range("B6").Value = $ws.Range("B5").Value()
$ws.Range("C6").Value = $ws.Range("C5").Value()

$ws.Range("A5").Value = $ws.Range("A4").Value()
$ws.Range("B5").Value = $ws.Range("B4").Value()
$ws.Range("C5").Value = $ws.Range("C4").Value()

$ws.Range("A4").Value = $ws.Range("A3").Value()
$ws.Range("B4").Value = $ws.Range("B3").Value()
$ws.Range("C4").Value = $ws.Range("C3").Value()

# New accessibility-check row at row 3.
$ws.Range("A3").Value = "checkAccessibility"
$ws.Range("B3").Value = "TC_PM_COA_SEC_ListView_D2"
$ws.Range("C3").ClearContents() | Out-Null

# B3 gets its own wrapped/top-left-aligned style (new cellXfs entry), based
# on the borderless/no-fill formatting already used by the F column.
$ws.Range("F2").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4122) | Out-Null
$ws.Range("B3").HorizontalAlignment = -4131
$ws.Range("B3").VerticalAlignment = -4160
$ws.Range("B3").WrapText = $true

# New trailing blank row (row 9), matching the formatting already used
# by row 8's D:F cells.
$ws.Range("D8:F8").Copy() | Out-Null
$ws.Range("D9:F9").PasteSpecial(-4122) | Out-Null
$ws.Range("D9:F9").ClearContents() | Out-Null

$ws.Range("C9").Select() | Out-Null

$ws2 = $wb.Worksheets.Item("TestData")
$ws2.Activate() | Out-Null
$ws2.Range("D12").Select() | Out-Null
